# Update cfb_weather.xlsx with Timestamp 2024-10-30T16:21:44.228691
# Applies the weather/odds-model refresh recorded in the commit diff:
#   - FBS sheet: assorted odds / wind / temp cell updates on several game rows
#   - FBS sheet: Timestamp column (AK) refreshed for every game row
#   - Other sheet: wind-direction / temp_fg updates on two rows
# Cell coordinates and values below were taken directly from the unified
# diff of the workbook's OOXML.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "FBS"
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("FBS")

$ws.Range("Y4").Value = 58.5
$ws.Range("AE4").Value = 0.01739130434782609

$ws.Range("Y5").Value = 45.5
$ws.Range("AE5").Value = -0.02150537634408602

$ws.Range("Y6").Value = 50.5
$ws.Range("AE6").Value = -0.01941747572815534

$ws.Range("Y7").Value = 49.5
$ws.Range("Z7").Value = -110
$ws.Range("AE7").Value = -0.05714285714285714

$ws.Range("AB8").Value = -13.5
$ws.Range("AF8").Value = -1.5

$ws.Range("Q9").Value = "NNW"

$ws.Range("Z10").Value = -105

$ws.Range("Y13").Value = 57.5
$ws.Range("Z13").Value = -106
$ws.Range("AE13").Value = 0

$ws.Range("O19").Value = 67.10000000000001
$ws.Range("P19").Value = 8.1
$ws.Range("U19").Value = 4.3
$ws.Range("Y19").Value = 64.5
$ws.Range("Z19").Value = -112
$ws.Range("AB19").Value = -2.5
$ws.Range("AE19").Value = 0.01574803149606299
$ws.Range("AF19").Value = 1

$ws.Range("O20").Value = 70.73
$ws.Range("P20").Value = 9.1
$ws.Range("U20").Value = 2.3

$ws.Range("O21").Value = 67.04000000000001
$ws.Range("P21").Value = 5.1
$ws.Range("U21").Value = -0.7
$ws.Range("AB21").Value = 14.5
$ws.Range("AF21").Value = 1

$ws.Range("Y23").Value = 42.5
$ws.Range("AB23").Value = -21.5
$ws.Range("AE23").Value = 0.04938271604938271
$ws.Range("AF23").Value = -2

$ws.Range("Z26").Value = -105

$ws.Range("Z35").Value = -118

$ws.Range("Y36").Value = 60.5
$ws.Range("Z36").Value = -114
$ws.Range("AE36").Value = -0.01626016260162602

$ws.Range("Q38").Value = "WNW"

$ws.Range("Z39").Value = -115

$ws.Range("AB40").Value = -16.5
$ws.Range("AF40").Value = 1

# Timestamp column: every game row (2-45) gets the refreshed run timestamp.
$newTimestamp = "2024-10-30T16:21:44.228691"
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 37).Value = $newTimestamp
}

# ----------------------------------------------------------------------
# Sheet "Other"
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Other")

$ws2.Range("O20").Value = "NNW"
$ws2.Range("P20").Value = "NW"
$ws2.Range("Q20").Value = 70.99999999999999
$ws2.Range("R20").Value = 6.4
$ws2.Range("S20").Value = "NNW"

$ws2.Range("S24").Value = "S"
